$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: "information card"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("information card")

# Row 9 becomes what used to be row 10 ("Action Resume" entry),
# gaining a "common" marker in column A.
$ws2.Range("A9").Value = "공통"
$ws2.Range("B9").Value = "38"
$ws2.Range("C9").Value = "행동 재개"
$ws2.Range("D9").Value = "특별"
$ws2.Range("G9").Value = "두 번 더 행동할 수 있습니다"

# Old row 10 is fully removed (not merely blanked).
$ws2.Range("A10:H10").Clear()

# Row 18 (previously an empty placeholder row) becomes what used to be row 9
# ("입막음" entry).
$ws2.Range("B18").NumberFormat = "@"
$ws2.Range("B18").HorizontalAlignment = -4108
$ws2.Range("B18").VerticalAlignment = -4108
$ws2.Range("B18").Value = "37"
$ws2.Range("C18").Value = "입막음"
$ws2.Range("D18").Value = "특별"
$ws2.Range("E18").Value = 0.05
$ws2.Range("F18").Value = 1
$ws2.Range("G18").Value = "대상은 잠시동안 발언할 수 없습니다"

# ---------------------------------------------------------------------
# Sheet: "battle card"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("battle card")

# Row 12 becomes what used to be row 13 ("행동 재개" entry), gaining the
# "common" marker in column A.
$ws3.Range("A12").Value = "공통"
$ws3.Range("B12").Value = "70"
$ws3.Range("C12").Value = "행동 재개"
$ws3.Range("D12").Value = "특별"
$ws3.Range("I12").Value = "두 번 더 행동할 수 있습니다"

# Old row 13 is fully removed (not merely blanked).
$ws3.Range("A13:J13").Clear()

# New row 20 gets what used to be row 12 ("금제" entry).
$ws3.Range("B20").NumberFormat = "@"
$ws3.Range("B20").HorizontalAlignment = -4108
$ws3.Range("B20").VerticalAlignment = -4108
$ws3.Range("B20").Value = "69"
$ws3.Range("C20").Value = "금제"
$ws3.Range("D20").Value = "특별"
$ws3.Range("G20").Value = 0.05
$ws3.Range("H20").Value = 1
$ws3.Range("I20").Value = "대상은 잠시동안 행동할 수 없습니다"

# battle card keeps its own remembered selection, but is no longer the
# active/front-most sheet.
$ws3.Range("E7").Select()

# ---------------------------------------------------------------------
# Make "information card" the active sheet/tab with its new selection.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F19").Select()
